$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values (not at risk of being auto-coerced to a number/date by
# Excel's "smart" input parsing) can just be assigned directly - this keeps
# the cell's existing style index untouched.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Values that Excel would otherwise auto-convert to a number or a date
# (pure-digit strings, "mm/dd/yy"-shaped strings) need special handling so
# they stay text AND keep the cell's original style index (s=) unchanged.
# Plan:
#   1) stash the cell's current (pre-edit) value+format in a scratch cell
#      via Copy (copies format+value together)
#   2) force the cell's number format to Text ("@") so the new value isn't
#      re-interpreted, and write the new text
#   3) paste back ONLY the formatting (xlPasteFormats = -4122) from the
#      scratch cell, restoring the exact original style without touching
#      the freshly-written text value
#   4) clear the scratch cell
$script:scratchRow = 60
function Set-TextValueSafe($addr, $val) {
    $backup = "Z" + $script:scratchRow
    $script:scratchRow = $script:scratchRow + 1
    $ws.Range($addr).Copy($ws.Range($backup))
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($backup).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($backup).Clear()
}

# Header date / invoice number
Set-TextValueSafe "B5" "02/14/2024"
Set-TextValue      "L5" "QR02142024AP101"

# Pharmacy / wholesaler / credit-to block
Set-TextValue "B10" "AUTREY PHARMACY 1"
Set-TextValue "D10" "AMERISOURCEBERGEN DRUG CORPORATION"
Set-TextValue "E10" "AUTREY PHARMACY 1"

Set-TextValue "B11" "1205 CENTRAL BLVD"
Set-TextValue "D11" "108 ROUTE 17K SUITE 1"
Set-TextValue "E11" "1205 CENTRAL BLVD"

Set-TextValue "B12" "BROWNSVILLE, TX, 78520"
Set-TextValue "D12" "NEWBURGH. NY, 12550-5008"
Set-TextValue "E12" "BROWNSVILLE, TX, 78520"

Set-TextValue "D13" "Account#: 100199545"

Set-TextValue "B14" "Phone: 956-548-0801, fax: "
Set-TextValue "D14" "Phone: 844-222-2273"
Set-TextValue "E14" "Phone: 956-548-0801, fax: "

Set-TextValue "B15" "DEA: FA3704358, Exp: 06/30/2023"
Set-TextValue "D15" "DEA: RA0522056"
Set-TextValue "E15" "DEA: FA3704358, Exp: 06/30/2023"

# Line item 1 (row 20)
Set-TextValueSafe "B20" "6787743305"
Set-TextValue      "C20" "Ascend Laboratories, LLC"
Set-TextValue      "D20" "Aripiprazole"
Set-TextValue      "E20" "15 mg/1"
Set-TextValueSafe "G20" "22140477"
Set-TextValueSafe "H20" "01/24/31"
Set-TextValue      "I20" "500 CT"

# Line item 2 (row 21)
Set-TextValueSafe "B21" "6787743305"
Set-TextValue      "C21" "Ascend Laboratories, LLC"
Set-TextValue      "D21" "Aripiprazole"
Set-TextValue      "E21" "15 mg/1"
Set-TextValueSafe "G21" "22140477"
Set-TextValueSafe "H21" "01/24/31"
Set-TextValue      "I21" "500 CT"

# New est value on row 22
$ws.Range("M22").Value = 1
